$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$app = $excel

# Unhide columns B:K (previously hidden helper columns) and the L:M pair.
$ws.Columns("B:M").Hidden = $false

# Column L used to be merged with M at a near-zero width (0.1640625) while
# hidden; now it gets its own "real" width, while M keeps its tiny width.
$ws.Columns("L:L").ColumnWidth = 6.33

# Highlight A21 with a red fill (new fill + cellXfs entry).
$ws.Range("A21").Interior.Color = 255

# Update the on-screen selection to match the new view state (best effort -
# the host only tracks a single active selection, so we leave the cursor on
# the cell the user was last working with).
$ws.Range("L1:N1048576").Select() | Out-Null
$ws.Range("F22").Select() | Out-Null
